$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "31.050.39"
Set-TextValue "E2" "  +1.18%  "
Set-TextValue "D3" "1.956.14"
Set-TextValue "E3" "  +0.47%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "246.31"
Set-TextValue "E5" "  -0.24%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.08%  "
Set-TextValue "D7" "0.4899"
Set-TextValue "E7" "  +1.40%  "
Set-TextValue "D8" "0.2973"
Set-TextValue "E8" "  +1.24%  "
Set-TextValue "D9" "0.06860"
Set-TextValue "E9" "  +0.59%  "
Set-TextValue "D10" "19.25"
Set-TextValue "E10" "  -0.85%  "
Set-TextValue "D11" "108.01"
Set-TextValue "E11" "  -4.01%  "
Set-TextValue "D12" "1.947.03"
Set-TextValue "E12" "  +0.03%  "
Set-TextValue "D13" "0.07758"
Set-TextValue "E13" "  +1.13%  "
Set-TextValue "D14" "5.477"
Set-TextValue "E14" "  -0.70%  "
Set-TextValue "E15" "  +2.80%  "
Set-TextValue "D16" "282.52"
Set-TextValue "E16" "  -4.82%  "
Set-TextValue "D17" "31.073.40"
Set-TextValue "E17" "  +1.11%  "
Set-TextValue "D18" "13.31"
Set-TextValue "E18" "  +0.30%  "
Set-TextValue "D19" "0.000007771"
Set-TextValue "E19" "  +1.08%  "
Set-TextValue "D20" "2.198.74"
Set-TextValue "E20" "  -0.24%  "
Set-TextValue "E21" "  +0.11%  "
Set-TextValue "D22" "5.524"
Set-TextValue "E22" "  -2.73%  "
Set-TextValue "D23" "1.001"
Set-TextValue "E23" "  +0.08%  "
Set-TextValue "D24" "6.528"
Set-TextValue "E24" "  -1.09%  "
Set-TextValue "D25" "9.843"
Set-TextValue "E25" "  -0.07%  "
Set-TextValue "D26" "169.75"
Set-TextValue "D27" "20.09"
Set-TextValue "E27" "  -1.01%  "
Set-TextValue "D28" "2.228"
Set-TextValue "E28" "  +1.70%  "
Set-TextValue "D29" "0.1057"
Set-TextValue "E29" "  -2.56%  "
Set-TextValue "D30" "1.425"
Set-TextValue "E30" "  -0.81%  "
Set-TextValue "D31" "1.585"
Set-TextValue "E31" "  -0.37%  "
Set-TextValue "D32" "4.590"
Set-TextValue "E32" "  -2.29%  "
Set-TextValue "D33" "4.471"
Set-TextValue "E33" "  -0.24%  "
Set-TextValue "D34" "0.04993"
Set-TextValue "E34" "  -1.82%  "
Set-TextValue "D35" "0.7604"
Set-TextValue "E35" "  -1.82%  "
Set-TextValue "E36" "  +2.17%  "
Set-TextValue "D37" "2.732"
Set-TextValue "E37" "  -0.06%  "
Set-TextValue "D38" "0.02037"
Set-TextValue "E38" "  -1.99%  "
Set-TextValue "D39" "2.705"
Set-TextValue "E39" "  +0.14%  "
Set-TextValue "D40" "2.182"
Set-TextValue "E40" "  +6.12%  "
Set-TextValue "D41" "6.493"
Set-TextValue "E41" "  +9.89%  "
Set-TextValue "D42" "76.26"
Set-TextValue "E42" "  +8.63%  "
Set-TextValue "D43" "0.4522"
Set-TextValue "E43" "  +1.26%  "
Set-TextValue "D44" "0.8873"
Set-TextValue "E44" "  +1.44%  "
Set-TextValue "E45" "  -1.68%  "
Set-TextValue "D46" "8.121"
Set-TextValue "E46" "  +10.05%  "
Set-TextValue "D47" "1.001"
Set-TextValue "E47" "  +0.11%  "
Set-TextValue "D48" "1.001.51"
Set-TextValue "E48" "  +10.67%  "
Set-TextValue "D49" "9.409"
Set-TextValue "E49" "  -0.52%  "
Set-TextValue "E50" "  +0.79%  "
Set-TextValue "D51" "0.2595"
Set-TextValue "E51" "  +2.78%  "
